$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (column D) and Volume(1h) (column E) columns with refreshed
# crypto data. All cells in these columns are stored as text, so values that
# look like plain numbers (e.g. "579.52") are written with a leading
# apostrophe to force Excel to keep them as text, then the cell style is
# reset to "Normal" so no extra number-format/quote-prefix styling leaks
# into the saved cell (matching the original file's plain default style).

$ws.Range("D2").Value = '65.030.48'
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").Value = '3.175.97'
$ws.Range("E3").Value = '  +4.02%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''579.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.65%  '
$ws.Range("D6").Value = '''151.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.14%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.175.97'
$ws.Range("E8").Value = '  +4.05%  '
$ws.Range("D9").Value = '''0.536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.28%  '
$ws.Range("E10").Value = '  +5.79%  '
$ws.Range("D11").Value = '''6.23'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  +18.59%  '
$ws.Range("D14").Value = '''37.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.15%  '
$ws.Range("D15").Value = '3.697.13'
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").Value = '65.122.92'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").Value = '3.177.80'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("E18").Value = '  +5.73%  '
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '''512.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.79%  '
$ws.Range("E21").Value = '  +6.35%  '
$ws.Range("D22").Value = '''0.731'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.94%  '
$ws.Range("D23").Value = '''15.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.21%  '
$ws.Range("E24").Value = '  +3.44%  '
$ws.Range("D25").Value = '''85.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  +10.90%  '
$ws.Range("D28").Value = '''2.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.65%  '
$ws.Range("E29").Value = '  +7.53%  '
$ws.Range("D30").Value = '''27.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("E31").Value = '  +13.73%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  +5.53%  '
$ws.Range("E34").Value = '  +10.22%  '
$ws.Range("D35").Value = '''6.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("D36").Value = '''55.72'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.23%  '
$ws.Range("D37").Value = '''0.0904'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.10%  '
$ws.Range("D38").Value = '''3.15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.32%  '
$ws.Range("D39").Value = '''474.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.22%  '
$ws.Range("E40").Value = '  +3.06%  '
$ws.Range("E41").Value = '  +4.75%  '
$ws.Range("D42").Value = '3.068.16'
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("E44").Value = '  +6.64%  '
$ws.Range("D45").Value = '''2.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.85%  '
$ws.Range("D46").Value = '''29.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.57%  '
$ws.Range("D47").Value = '0.0₃0614'
$ws.Range("E47").Value = '  +19.53%  '
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("D50").Value = '''2.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.35%  '
$ws.Range("D51").Value = '''120.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.59%  '
